$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old Transaction ID data (rows 1-22) and replace with the
# new two-row Corporate Customer data set.
$ws.Range("A1:A22").ClearContents()

$ws.Range("A1").Value = "Transaction Number"
$ws.Range("A2").Value = "FT231850FF7GLCWC"

# Match the new selected cell recorded in the saved worksheet view.
$ws.Range("L6").Select()
